$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 9114
$ws.Range("I18").Value = 8966.333000000001
$ws.Range("K18").Value = 8966.333000000001
$ws.Range("M18").Value = -8682.333000000001
$ws.Range("H28").Value = 1243.25
$ws.Range("I28").Value = 1424.6
$ws.Range("J28").Value = 336.5
$ws.Range("K28").Value = 1424.6
$ws.Range("L28").Value = 336.5
$ws.Range("M28").Value = -939.5999999999999
$ws.Range("N28").Value = -1306.5
$ws.Range("H32").Value = 2724.8462
$ws.Range("J32").Value = 3956.7144
$ws.Range("L32").Value = 3956.7144
$ws.Range("N32").Value = -4608.7144
$ws.Range("H80").Value = 2773.25
$ws.Range("I80").Value = 2666.6667
$ws.Range("J80").Value = 2808.7778
$ws.Range("K80").Value = 8000.000100000001
$ws.Range("L80").Value = 8426.3334
$ws.Range("M80").Value = -7002.000100000001
$ws.Range("N80").Value = -10422.3334
$ws.Range("H83").Value = 2773.25
$ws.Range("I83").Value = 2666.6667
$ws.Range("J83").Value = 2808.7778
$ws.Range("K83").Value = 24000.0003
$ws.Range("L83").Value = 25279.0002
$ws.Range("M83").Value = -19008.0003
$ws.Range("N83").Value = -35263.00019999999
$ws.Range("H88").Value = 1699.3334
$ws.Range("J88").Value = 1699.3334
$ws.Range("L88").Value = 1699.3334
$ws.Range("N88").Value = -2511.3334
$ws.Range("H91").Value = 1699.3334
$ws.Range("J91").Value = 1699.3334
$ws.Range("L91").Value = 1699.3334
$ws.Range("N91").Value = -4507.3334
$ws.Range("H103").Value = 668.25
$ws.Range("J103").Value = 799.4
$ws.Range("L103").Value = 2398.2
$ws.Range("N103").Value = -3570.2
$ws.Range("H138").Value = 3998.647
$ws.Range("I138").Value = 2359.818
$ws.Range("J138").Value = 4782.4346
$ws.Range("K138").Value = 7079.454000000001
$ws.Range("L138").Value = 14347.3038
$ws.Range("M138").Value = -1939.454000000001
$ws.Range("N138").Value = -24627.3038
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3329.3333
$ws.Range("I2").Value = 2990
$ws.Range("K2").Value = 2990
$ws.Range("M2").Value = -2877
$ws.Range("H32").Value = 3211785
$ws.Range("I32").Value = 5389314
$ws.Range("J32").Value = 638341.4
$ws.Range("K32").Value = 5389314
$ws.Range("L32").Value = 638341.4
$ws.Range("M32").Value = -5389027
$ws.Range("N32").Value = -638915.4
$ws.Range("H61").Value = 3258.1667
$ws.Range("I61").Value = 3227.0908
$ws.Range("K61").Value = 3227.0908
$ws.Range("M61").Value = -3015.0908
$ws.Range("H116").Value = 3329.3333
$ws.Range("I116").Value = 2990
$ws.Range("K116").Value = 2990
$ws.Range("M116").Value = -696
$ws.Range("H132").Value = 2445.6667
$ws.Range("I132").Value = 2445.6667
$ws.Range("K132").Value = 7337.000100000001
$ws.Range("M132").Value = -4807.000100000001
$ws.Range("H136").Value = 3258.1667
$ws.Range("I136").Value = 3227.0908
$ws.Range("K136").Value = 9681.2724
$ws.Range("M136").Value = -7131.2724
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3329.3333
$ws.Range("I3").Value = 2990
$ws.Range("K3").Value = 2990
$ws.Range("M3").Value = -2876
$ws.Range("H20").Value = 1950.2
$ws.Range("I20").Value = 1139.3334
$ws.Range("J20").Value = 3166.5
$ws.Range("K20").Value = 1139.3334
$ws.Range("L20").Value = 3166.5
$ws.Range("M20").Value = -892.3334
$ws.Range("N20").Value = -3660.5
$ws.Range("H76").Value = 144841.67
$ws.Range("J76").Value = 144841.67
$ws.Range("L76").Value = 144841.67
$ws.Range("N76").Value = -145471.67
$ws.Range("H79").Value = 144841.67
$ws.Range("J79").Value = 144841.67
$ws.Range("L79").Value = 144841.67
$ws.Range("N79").Value = -147025.67
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 114.833336
$ws.Range("I7").Value = 94
$ws.Range("J7").Value = 119
$ws.Range("K7").Value = 94
$ws.Range("L7").Value = 119
$ws.Range("M7").Value = 19
$ws.Range("N7").Value = -345
$ws.Range("H22").Value = 494.75
$ws.Range("I22").Value = 494.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 494.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -144.75
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 2214.0908
$ws.Range("I58").Value = 1990.1111
$ws.Range("K58").Value = 1990.1111
$ws.Range("M58").Value = -1787.1111
$ws.Range("H86").Value = 19165.166
$ws.Range("I86").Value = 19997.75
$ws.Range("K86").Value = 19997.75
$ws.Range("M86").Value = -18874.75
$ws.Range("H89").Value = 19165.166
$ws.Range("I89").Value = 19997.75
$ws.Range("K89").Value = 99988.75
$ws.Range("M89").Value = -94372.75
$ws.Range("H107").Value = 363.5
$ws.Range("I107").Value = 332.33334
$ws.Range("J107").Value = 394.66666
$ws.Range("K107").Value = 332.33334
$ws.Range("L107").Value = 394.66666
$ws.Range("M107").Value = 1587.66666
$ws.Range("N107").Value = -4234.66666
$ws.Range("H132").Value = 4030.8333
$ws.Range("I132").Value = 4030.8333
$ws.Range("K132").Value = 12092.4999
$ws.Range("M132").Value = -9562.499899999999
$ws.Range("H136").Value = 2214.0908
$ws.Range("I136").Value = 1990.1111
$ws.Range("K136").Value = 5970.3333
$ws.Range("M136").Value = -3420.3333
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 5
$ws.Range("K19").Value = 5
$ws.Range("M19").Value = 283
$ws.Range("H70").Value = 6400
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 6400
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H122").Value = 1980.2858
$ws.Range("I122").Value = 1994
$ws.Range("J122").Value = 1898
$ws.Range("K122").Value = 5982
$ws.Range("L122").Value = 5694
$ws.Range("M122").Value = -3532
$ws.Range("N122").Value = -10594
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 1250
$ws.Range("I24").Value = 1250
$ws.Range("K24").Value = 1250
$ws.Range("M24").Value = -907
$ws.Range("H68").Value = 549.5
$ws.Range("I68").Value = 549.5
$ws.Range("K68").Value = 549.5
$ws.Range("M68").Value = 199.5
$ws.Range("H69").Value = 35494.5
$ws.Range("J69").Value = 35494.5
$ws.Range("L69").Value = 35494.5
$ws.Range("N69").Value = -37116.5
$ws.Range("H71").Value = 549.5
$ws.Range("I71").Value = 549.5
$ws.Range("K71").Value = 2747.5
$ws.Range("M71").Value = 996.5
$ws.Range("H72").Value = 35494.5
$ws.Range("J72").Value = 35494.5
$ws.Range("L72").Value = 106483.5
$ws.Range("N72").Value = -114595.5
$ws.Range("I122").Value = 4908.1665
$ws.Range("J122").Value = 7262.5713
$ws.Range("K122").Value = 14724.4995
$ws.Range("L122").Value = 21787.7139
$ws.Range("M122").Value = -12274.4995
$ws.Range("N122").Value = -26687.7139
$ws.Range("H132").Value = 8908.777
$ws.Range("I132").Value = 8908.777
$ws.Range("K132").Value = 26726.331
$ws.Range("M132").Value = -24196.331
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 26840.8
$ws.Range("I62").Value = 55002
$ws.Range("J62").Value = 8066.6665
$ws.Range("K62").Value = 55002
$ws.Range("L62").Value = 8066.6665
$ws.Range("M62").Value = -54378
$ws.Range("N62").Value = -9314.666499999999
$ws.Range("H63").Value = 41666
$ws.Range("J63").Value = 41666
$ws.Range("L63").Value = 41666
$ws.Range("N63").Value = -42914
$ws.Range("H65").Value = 26840.8
$ws.Range("I65").Value = 55002
$ws.Range("J65").Value = 8066.6665
$ws.Range("K65").Value = 275010
$ws.Range("L65").Value = 40333.3325
$ws.Range("M65").Value = -271890
$ws.Range("N65").Value = -46573.3325
$ws.Range("H66").Value = 41666
$ws.Range("J66").Value = 41666
$ws.Range("L66").Value = 124998
$ws.Range("N66").Value = -131238
$ws.Range("H68").Value = 43300
$ws.Range("I68").Value = 25000
$ws.Range("J68").Value = 52450
$ws.Range("K68").Value = 25000
$ws.Range("L68").Value = 52450
$ws.Range("M68").Value = -24189
$ws.Range("N68").Value = -54072
$ws.Range("H71").Value = 43300
$ws.Range("I71").Value = 25000
$ws.Range("J71").Value = 52450
$ws.Range("K71").Value = 75000
$ws.Range("L71").Value = 157350
$ws.Range("M71").Value = -70944
$ws.Range("N71").Value = -165462
$ws.Range("H81").Value = 911225.6
$ws.Range("I81").Value = 1282.2858
$ws.Range("J81").Value = 2503626.5
$ws.Range("K81").Value = 2564.5716
$ws.Range("L81").Value = 5007253
$ws.Range("M81").Value = -1503.5716
$ws.Range("N81").Value = -5009375
$ws.Range("H84").Value = 911225.6
$ws.Range("I84").Value = 1282.2858
$ws.Range("J84").Value = 2503626.5
$ws.Range("K84").Value = 12822.858
$ws.Range("L84").Value = 25036265
$ws.Range("M84").Value = -7518.858
$ws.Range("N84").Value = -25046873
$ws.Range("H136").Value = 2956.1304
$ws.Range("I136").Value = 3085.3809
$ws.Range("K136").Value = 9256.1427
$ws.Range("M136").Value = -6706.1427
